# Auto-generated edit: refresh Market Board price columns (H-N) on several
# "Leve Profits" rows across the ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Columns: H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#          K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ARM!row110 (Leve Item ID 27708)
$ws_ARM.Range("H110").Value = 762.8570999999999  # was 637.05884
$ws_ARM.Range("I110").Value = 640.0833  # was 547.53845
$ws_ARM.Range("J110").Value = 1499.5  # was 928
$ws_ARM.Range("K110").Value = 640.0833  # was 547.53845
$ws_ARM.Range("L110").Value = 1499.5  # was 928
$ws_ARM.Range("M110").Value = 1404.9167  # was 1497.46155
$ws_ARM.Range("N110").Value = -5589.5  # was -5018

# BSM!row20 (Leve Item ID 14149)
$ws_BSM.Range("H20").Value = 2918  # was 3059.9333
$ws_BSM.Range("I20").Value = 3049.2144  # was 3223.077
$ws_BSM.Range("K20").Value = 3049.2144  # was 3223.077
$ws_BSM.Range("M20").Value = -2802.2144  # was -2976.077

# BSM!row86 (Leve Item ID 12526)
$ws_BSM.Range("H86").Value = 1570.9811  # was 1495.9822
$ws_BSM.Range("I86").Value = 1378.0264  # was 1329.8718
$ws_BSM.Range("J86").Value = 2059.8  # was 1877.0588
$ws_BSM.Range("K86").Value = 1378.0264  # was 1329.8718
$ws_BSM.Range("L86").Value = 2059.8  # was 1877.0588
$ws_BSM.Range("M86").Value = -255.0264  # was -206.8717999999999
$ws_BSM.Range("N86").Value = -4305.8  # was -4123.0588

# BSM!row89 (Leve Item ID 12526)
$ws_BSM.Range("H89").Value = 1570.9811  # was 1495.9822
$ws_BSM.Range("I89").Value = 1378.0264  # was 1329.8718
$ws_BSM.Range("J89").Value = 2059.8  # was 1877.0588
$ws_BSM.Range("K89").Value = 6890.132  # was 6649.358999999999
$ws_BSM.Range("L89").Value = 10299  # was 9385.294
$ws_BSM.Range("M89").Value = -1274.132  # was -1033.358999999999
$ws_BSM.Range("N89").Value = -21531  # was -20617.294

# BSM!row105 (Leve Item ID 19947)
$ws_BSM.Range("H105").Value = 3800  # was 3475
$ws_BSM.Range("I105").Value = 3150  # was 2766.6667
$ws_BSM.Range("J105").Value = 4125  # was 5600
$ws_BSM.Range("K105").Value = 3150  # was 2766.6667
$ws_BSM.Range("L105").Value = 4125  # was 5600
$ws_BSM.Range("M105").Value = -1403  # was -1019.6667
$ws_BSM.Range("N105").Value = -7619  # was -9094

# BSM!row107 (Leve Item ID 27706)
$ws_BSM.Range("H107").Value = 1406.8  # was 1592.4
$ws_BSM.Range("I107").Value = 756.875  # was 801.375
$ws_BSM.Range("J107").Value = 4006.5  # was 4756.5
$ws_BSM.Range("K107").Value = 756.875  # was 801.375
$ws_BSM.Range("L107").Value = 4006.5  # was 4756.5
$ws_BSM.Range("M107").Value = 1163.125  # was 1118.625
$ws_BSM.Range("N107").Value = -7846.5  # was -8596.5

# BSM!row131 (Leve Item ID 35396)
$ws_BSM.Range("H131").Value = 24997.846  # was 26185.188
$ws_BSM.Range("J131").Value = 24997.846  # was 26185.188
$ws_BSM.Range("L131").Value = 24997.846  # was 26185.188
$ws_BSM.Range("N131").Value = -35077.84600000001  # was -36265.18799999999

# CRP!row41 (Leve Item ID 1917)
$ws_CRP.Range("H41").Value = 0  # was 5059
$ws_CRP.Range("I41").Value = 0  # was 5059
$ws_CRP.Range("K41").Value = 0  # was 5059
$ws_CRP.Range("M41").ClearContents()  # was -4631, now blank

# CRP!row58 (Leve Item ID 44021)
$ws_CRP.Range("H58").Value = 20434  # was 19701.26
$ws_CRP.Range("I58").Value = 1327.8235  # was 1290.1666
$ws_CRP.Range("K58").Value = 1327.8235  # was 1290.1666
$ws_CRP.Range("M58").Value = -1124.8235  # was -1087.1666

# CRP!row59 (Leve Item ID 1942)
$ws_CRP.Range("H59").Value = 20647  # was 20163.316
$ws_CRP.Range("I59").Value = 4999  # was 7551.5
$ws_CRP.Range("J59").Value = 21625  # was 21647.059
$ws_CRP.Range("K59").Value = 4999  # was 7551.5
$ws_CRP.Range("L59").Value = 21625  # was 21647.059
$ws_CRP.Range("M59").Value = -3854  # was -6406.5
$ws_CRP.Range("N59").Value = -23915  # was -23937.059

# CRP!row60 (Leve Item ID 1937)
$ws_CRP.Range("H60").Value = 11135.412  # was 10173.685
$ws_CRP.Range("I60").Value = 5000  # was 2999.3333
$ws_CRP.Range("K60").Value = 5000  # was 2999.3333
$ws_CRP.Range("M60").Value = -4489  # was -2488.3333

# CRP!row105 (Leve Item ID 19928)
$ws_CRP.Range("H105").Value = 9616652  # was 9616668
$ws_CRP.Range("I105").Value = 15625688  # was 20834056
$ws_CRP.Range("J105").Value = 2194.4  # was 1762.8572
$ws_CRP.Range("K105").Value = 15625688  # was 20834056
$ws_CRP.Range("L105").Value = 2194.4  # was 1762.8572
$ws_CRP.Range("M105").Value = -15623941  # was -20832309
$ws_CRP.Range("N105").Value = -5688.4  # was -5256.8572

# CRP!row132 (Leve Item ID 44019)
$ws_CRP.Range("H132").Value = 15741.711  # was 14979.625
$ws_CRP.Range("I132").Value = 22017.44  # was 20423.555
$ws_CRP.Range("K132").Value = 66052.31999999999  # was 61270.665
$ws_CRP.Range("M132").Value = -63522.31999999999  # was -58740.665

# CRP!row134 (Leve Item ID 44020)
$ws_CRP.Range("H134").Value = 1103.5625  # was 1123.6774
$ws_CRP.Range("I134").Value = 1034.5834  # was 1085
$ws_CRP.Range("K134").Value = 3103.7502  # was 3255
$ws_CRP.Range("M134").Value = -568.7501999999999  # was -720

# CRP!row136 (Leve Item ID 44021)
$ws_CRP.Range("H136").Value = 20434  # was 19701.26
$ws_CRP.Range("I136").Value = 1327.8235  # was 1290.1666
$ws_CRP.Range("K136").Value = 3983.4705  # was 3870.4998
$ws_CRP.Range("M136").Value = -1433.4705  # was -1320.4998

# CUL!row98 (Leve Item ID 19843)
$ws_CUL.Range("H98").Value = 583.8333  # was 349.92307
$ws_CUL.Range("I98").Value = 499  # was 256.75
$ws_CUL.Range("J98").Value = 600.8  # was 391.33334
$ws_CUL.Range("K98").Value = 1497  # was 770.25
$ws_CUL.Range("L98").Value = 1802.4  # was 1174.00002
$ws_CUL.Range("M98").Value = 1  # was 727.75
$ws_CUL.Range("N98").Value = -4798.4  # was -4170.000019999999

# CUL!row131 (Leve Item ID 36060)
$ws_CUL.Range("H131").Value = 728.7778  # was 727.75
$ws_CUL.Range("I131").Value = 700  # was 0
$ws_CUL.Range("J131").Value = 729.0714  # was 727.75
$ws_CUL.Range("K131").Value = 2100  # was 0
$ws_CUL.Range("L131").Value = 2187.2142  # was 2183.25
$ws_CUL.Range("M131").Value = 2940  # newly added cell
$ws_CUL.Range("N131").Value = -12267.2142  # was -12263.25

# CUL!row140 (Leve Item ID 44097)
$ws_CUL.Range("H140").Value = 1633.421  # was 1564.2
$ws_CUL.Range("I140").Value = 1441.1765  # was 1374.9445
$ws_CUL.Range("K140").Value = 4323.529500000001  # was 4124.833500000001
$ws_CUL.Range("M140").Value = 856.4704999999994  # was 1055.166499999999

# GSM!row122 (Leve Item ID 36182)
$ws_GSM.Range("H122").Value = 4000  # was 1758.6
$ws_GSM.Range("I122").Value = 4000  # was 1509.5555
$ws_GSM.Range("K122").Value = 12000  # was 4528.666499999999
$ws_GSM.Range("M122").Value = -9550  # was -2078.666499999999

# GSM!row131 (Leve Item ID 34747)
$ws_GSM.Range("H131").Value = 40412.25  # was 37663
$ws_GSM.Range("J131").Value = 40412.25  # was 37663
$ws_GSM.Range("L131").Value = 40412.25  # was 37663
$ws_GSM.Range("N131").Value = -50492.25  # was -47743

# GSM!row132 (Leve Item ID 44008)
$ws_GSM.Range("H132").Value = 57871.535  # was 59948.63
$ws_GSM.Range("I132").Value = 60533.055  # was 54660.25
$ws_GSM.Range("J132").Value = 53080.8  # was 75058.28999999999
$ws_GSM.Range("K132").Value = 181599.165  # was 163980.75
$ws_GSM.Range("L132").Value = 159242.4  # was 225174.87
$ws_GSM.Range("M132").Value = -179069.165  # was -161450.75
$ws_GSM.Range("N132").Value = -164302.4  # was -230234.87

# LTW!row7 (Leve Item ID 36249)
$ws_LTW.Range("H7").Value = 3070.25  # was 3213.0715
$ws_LTW.Range("I7").Value = 3213.7896  # was 3413.7368
$ws_LTW.Range("J7").Value = 2767.2222  # was 2789.4443
$ws_LTW.Range("K7").Value = 3213.7896  # was 3413.7368
$ws_LTW.Range("L7").Value = 2767.2222  # was 2789.4443
$ws_LTW.Range("M7").Value = -3101.7896  # was -3301.7368
$ws_LTW.Range("N7").Value = -2991.2222  # was -3013.4443

# LTW!row22 (Leve Item ID 5277)
$ws_LTW.Range("H22").Value = 3540.1  # was 3381.9092
$ws_LTW.Range("I22").Value = 3180.2  # was 3140.2
$ws_LTW.Range("J22").Value = 3900  # was 3583.3333
$ws_LTW.Range("K22").Value = 3180.2  # was 3140.2
$ws_LTW.Range("L22").Value = 3900  # was 3583.3333
$ws_LTW.Range("M22").Value = -2885.2  # was -2845.2
$ws_LTW.Range("N22").Value = -4490  # was -4173.3333

# LTW!row26 (Leve Item ID 3559)
$ws_LTW.Range("H26").Value = 3803.3333  # was 3940
$ws_LTW.Range("J26").Value = 3803.3333  # was 3940
$ws_LTW.Range("L26").Value = 3803.3333  # was 3940
$ws_LTW.Range("N26").Value = -4393.3333  # was -4530

# LTW!row27 (Leve Item ID 5277)
$ws_LTW.Range("H27").Value = 3540.1  # was 3381.9092
$ws_LTW.Range("I27").Value = 3180.2  # was 3140.2
$ws_LTW.Range("J27").Value = 3900  # was 3583.3333
$ws_LTW.Range("K27").Value = 3180.2  # was 3140.2
$ws_LTW.Range("L27").Value = 3900  # was 3583.3333
$ws_LTW.Range("M27").Value = -3073.2  # was -3033.2
$ws_LTW.Range("N27").Value = -4114  # was -3797.3333

# LTW!row46 (Leve Item ID 5282)
$ws_LTW.Range("H46").Value = 1155.9286  # was 1298.4445
$ws_LTW.Range("I46").Value = 958.3  # was 997.6667
$ws_LTW.Range("J46").Value = 1650  # was 1900
$ws_LTW.Range("K46").Value = 958.3  # was 997.6667
$ws_LTW.Range("L46").Value = 1650  # was 1900
$ws_LTW.Range("M46").Value = -770.3  # was -809.6667
$ws_LTW.Range("N46").Value = -2026  # was -2276

# LTW!row62 (Leve Item ID 10740)
$ws_LTW.Range("H62").Value = 8000  # was 0
$ws_LTW.Range("J62").Value = 8000  # was 0
$ws_LTW.Range("L62").Value = 8000  # was 0
$ws_LTW.Range("N62").Value = -9248  # newly added cell

# LTW!row65 (Leve Item ID 10740)
$ws_LTW.Range("H65").Value = 8000  # was 0
$ws_LTW.Range("J65").Value = 8000  # was 0
$ws_LTW.Range("L65").Value = 24000  # was 0
$ws_LTW.Range("N65").Value = -30240  # newly added cell

# LTW!row126 (Leve Item ID 36249)
$ws_LTW.Range("H126").Value = 3070.25  # was 3213.0715
$ws_LTW.Range("I126").Value = 3213.7896  # was 3413.7368
$ws_LTW.Range("J126").Value = 2767.2222  # was 2789.4443
$ws_LTW.Range("K126").Value = 9641.3688  # was 10241.2104
$ws_LTW.Range("L126").Value = 8301.6666  # was 8368.332900000001
$ws_LTW.Range("M126").Value = -7171.3688  # was -7771.2104
$ws_LTW.Range("N126").Value = -13241.6666  # was -13308.3329

# LTW!row136 (Leve Item ID 44060)
$ws_LTW.Range("H136").Value = 22065.666  # was 19679.297
$ws_LTW.Range("I136").Value = 32068.5  # was 27097.947
$ws_LTW.Range("K136").Value = 96205.5  # was 81293.841
$ws_LTW.Range("M136").Value = -93655.5  # was -78743.841

# WVR!row96 (Leve Item ID 19977)
$ws_WVR.Range("H96").Value = 880  # was 1428.5714
$ws_WVR.Range("I96").Value = 500  # was 1250
$ws_WVR.Range("J96").Value = 1450  # was 2500
$ws_WVR.Range("K96").Value = 500  # was 1250
$ws_WVR.Range("L96").Value = 1450  # was 2500
$ws_WVR.Range("M96").Value = 873  # was 123
$ws_WVR.Range("N96").Value = -4196  # was -5246

# WVR!row122 (Leve Item ID 36208)
$ws_WVR.Range("H122").Value = 1355.7727  # was 1420.9512
$ws_WVR.Range("I122").Value = 1340.909  # was 1450
$ws_WVR.Range("J122").Value = 1400.3636  # was 1350.75
$ws_WVR.Range("K122").Value = 4022.727  # was 4350
$ws_WVR.Range("L122").Value = 4201.0908  # was 4052.25
$ws_WVR.Range("M122").Value = -1572.727  # was -1900
$ws_WVR.Range("N122").Value = -9101.0908  # was -8952.25
